$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-36, replacing the old Strike# derived
# values with the regenerated K values (per commit: "use K instead of Strike#,
# regen std/mean, calc and write s_vals").
$gValues = @{
    2  = 2
    3  = 3
    4  = 1
    5  = 3
    6  = 3
    7  = 6
    8  = 5
    9  = 8
    10 = 5
    11 = 1
    12 = 10
    13 = 11
    14 = 3
    15 = 3
    16 = 1
    17 = 6
    18 = 4
    19 = 2
    20 = 6
    21 = 5
    22 = 4
    23 = 3
    24 = 3
    25 = 2
    26 = 3
    27 = 1
    28 = 3
    29 = 4
    30 = 2
    31 = 1
    32 = 2
    33 = 5
    34 = 2
    35 = 8
    36 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
